$d = $word.ActiveDocument

# ---- Phase 1: replace each original text with a unique placeholder token ----
# (avoids collisions since this is a cyclic rotation of content between sections)
$d.Content.Find.Execute("Fornecer aos alunos os conhecimentos básicos dos sistemas principais de tratamento de água para consumo humano.", $true, $false, $false, $false, $false, $true, 1, $false, "@@P6@@", 2) | Out-Null
$d.Content.Find.Execute("Supply the students the basic knowledge of the systems of main water treatment technologies for the human consumption.", $true, $false, $false, $false, $false, $true, 1, $false, "@@P7@@", 2) | Out-Null
$d.Content.Find.Execute("7455355 - Robson da Silva Rocha", $true, $false, $false, $false, $false, $true, 1, $false, "@@P9@@", 2) | Out-Null
$d.Content.Find.Execute("Tecnologias de Tratamento de Água; Tratamento de Água em Ciclo Completo; Desinfecção; Filtração Direta Ascendente; Filtração Direta Descendente; Floto-Filtração; Filtração em Múltiplas Etapas; Tratamento dos Resíduos Gerados nas Estações de Tratamento de Água.", $true, $false, $false, $false, $false, $true, 1, $false, "@@P11@@", 2) | Out-Null
$d.Content.Find.Execute("Water Treatment Technologies; Water Treatment in Complete Cycle; Disinfection, Direct Ascendant Filtration; Direct Descendant Filtration, Floto-filtration; Filtration in Multiple Levels, Treatment of the Generated Waste in the Water Treatment Stations.", $true, $false, $false, $false, $false, $true, 1, $false, "@@P12@@", 2) | Out-Null
$d.Content.Find.Execute("- Características das águas de interesse para o tratamento: características físicas, químicas e bacteriológicas; `v- Padrão de Potabilidade; `v- Tecnologias de Tratamento de Água;`v- Unidades Constituintes de um Sistema de Abastecimento de Água;`v- Captação de Águas Subterrâneas e Captação de Águas Superficiais `v- Gradeamento, remoção de areia, casa de bombas; `v- Reservação; `v- Sistema de Tratamento de Água de Ciclo Completo; `v- Coagulação-floculação e Mistura Rápida; `v- Decantação: decantação convencional e de alta taxa e sistema de remoção de lodo;`v- Mecanismos da filtração, materiais filtrantes e fundos de filtros, hidráulica da filtração, filtração com taxa constante e taxa declinante, `v- Desinfecção: principais desinfetantes, cloração e cloro-amoniação, pré e pós-cloração, `v- Tratamento de resíduos gerados nas ETAs e reuso de água recuperada `v- A disciplina pode contar com viagens didáticas para complementação do conteúdo da disciplina.", $true, $false, $false, $false, $false, $true, 1, $false, "@@P14@@", 2) | Out-Null
$d.Content.Find.Execute("PROGRAMA DE PESQUISA EM SANEAMENTO BÁSICO. Tratamento de Água de`vAbastecimento por Filtração em Múltiplas Etapas. ASSOCIAÇÃO BRASILEIRA DE`vENGENHARIA SANITÁRIA E AMBIENTAL, Rio de Janeiro, 1999. (Coordenação: Luiz Di`vBernardo)." + ([char]0x00A0) + "`vPROGRAMA DE PESQUISA EM SANEAMENTO BÁSICO. Noções Gerais de Tratamento e`vDisposição Final de Lodos e Estações de Tratamento de Água. ASSOCIAÇÃO`vBRASILEIRA DE ENGENHARIA SANITÁRIA E AMBIENTAL, Rio de Janeiro, 200`v(Coordenação: Marco A.P. Reali).", $true, $false, $false, $false, $false, $true, 1, $false, "@@P19@@", 2) | Out-Null
$d.Content.Find.Execute("Avaliação baseada em provas, exercícios, trabalhos práticos e relatórios.`v", $true, $false, $false, $false, $false, $true, 1, $false, "@@P17R1@@", 2) | Out-Null
$d.Content.Find.Execute("Média ponderada das notas atribuídas às provas, exercícios e trabalhos práticos e relatórios.`v", $true, $false, $false, $false, $false, $true, 1, $false, "@@P17R3@@", 2) | Out-Null
$d.Content.Find.Execute("1 (uma) prova de recuperação (R), sendo considerado aprovado se 0,5(NF + R) >= 5,0.", $true, $false, $false, $false, $false, $true, 1, $false, "@@P17R5@@", 2) | Out-Null

# ---- Phase 2: replace each placeholder with its final destination text ----
$d.Content.Find.Execute("@@P6@@", $true, $false, $false, $false, $false, $true, 1, $false, "Tecnologias de Tratamento de Água; Tratamento de Água em Ciclo Completo; Desinfecção; Filtração Direta Ascendente; Filtração Direta Descendente; Floto-Filtração; Filtração em Múltiplas Etapas; Tratamento dos Resíduos Gerados nas Estações de Tratamento de Água.", 2) | Out-Null
$d.Content.Find.Execute("@@P7@@", $true, $false, $false, $false, $false, $true, 1, $false, "Water Treatment Technologies; Water Treatment in Complete Cycle; Disinfection, Direct Ascendant Filtration; Direct Descendant Filtration, Floto-filtration; Filtration in Multiple Levels, Treatment of the Generated Waste in the Water Treatment Stations.", 2) | Out-Null
$d.Content.Find.Execute("@@P9@@", $true, $false, $false, $false, $false, $true, 1, $false, "Fornecer aos alunos os conhecimentos básicos dos sistemas principais de tratamento de água para consumo humano.", 2) | Out-Null
$d.Content.Find.Execute("@@P11@@", $true, $false, $false, $false, $false, $true, 1, $false, "- Características das águas de interesse para o tratamento: características físicas, químicas e bacteriológicas; `v- Padrão de Potabilidade; `v- Tecnologias de Tratamento de Água;`v- Unidades Constituintes de um Sistema de Abastecimento de Água;`v- Captação de Águas Subterrâneas e Captação de Águas Superficiais `v- Gradeamento, remoção de areia, casa de bombas; `v- Reservação; `v- Sistema de Tratamento de Água de Ciclo Completo; `v- Coagulação-floculação e Mistura Rápida; `v- Decantação: decantação convencional e de alta taxa e sistema de remoção de lodo;`v- Mecanismos da filtração, materiais filtrantes e fundos de filtros, hidráulica da filtração, filtração com taxa constante e taxa declinante, `v- Desinfecção: principais desinfetantes, cloração e cloro-amoniação, pré e pós-cloração, `v- Tratamento de resíduos gerados nas ETAs e reuso de água recuperada `v- A disciplina pode contar com viagens didáticas para complementação do conteúdo da disciplina.", 2) | Out-Null
$d.Content.Find.Execute("@@P12@@", $true, $false, $false, $false, $false, $true, 1, $false, "Supply the students the basic knowledge of the systems of main water treatment technologies for the human consumption.", 2) | Out-Null
$d.Content.Find.Execute("@@P14@@", $true, $false, $false, $false, $false, $true, 1, $false, "Avaliação baseada em provas, exercícios, trabalhos práticos e relatórios.", 2) | Out-Null
$d.Content.Find.Execute("@@P19@@", $true, $false, $false, $false, $false, $true, 1, $false, "7455355 - Robson da Silva Rocha", 2) | Out-Null
$d.Content.Find.Execute("@@P17R1@@", $true, $false, $false, $false, $false, $true, 1, $false, "Média ponderada das notas atribuídas às provas, exercícios e trabalhos práticos e relatórios.`v", 2) | Out-Null
$d.Content.Find.Execute("@@P17R3@@", $true, $false, $false, $false, $false, $true, 1, $false, "1 (uma) prova de recuperação (R), sendo considerado aprovado se 0,5(NF + R) >= 5,0.`v", 2) | Out-Null
$d.Content.Find.Execute("@@P17R5@@", $true, $false, $false, $false, $false, $true, 1, $false, "PROGRAMA DE PESQUISA EM SANEAMENTO BÁSICO. Tratamento de Água de`vAbastecimento por Filtração em Múltiplas Etapas. ASSOCIAÇÃO BRASILEIRA DE`vENGENHARIA SANITÁRIA E AMBIENTAL, Rio de Janeiro, 1999. (Coordenação: Luiz Di`vBernardo)." + ([char]0x00A0) + "`vPROGRAMA DE PESQUISA EM SANEAMENTO BÁSICO. Noções Gerais de Tratamento e`vDisposição Final de Lodos e Estações de Tratamento de Água. ASSOCIAÇÃO`vBRASILEIRA DE ENGENHARIA SANITÁRIA E AMBIENTAL, Rio de Janeiro, 200`v(Coordenação: Marco A.P. Reali).", 2) | Out-Null
